# Update the cryptocurrency price/volume data on the active worksheet.
# Mirrors the GitHub Actions bot commit that refreshes "cryptos.xlsx"
# with the latest price (column D) and 1-hour volume change (column E)
# figures for each coin row (rows 2-51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "37.043.53"; E = "  -0.25%  " },
    @{ Row = 3; D = "2.058.98"; E = "  +0.10%  " },
    @{ Row = 4; D = $null; E = "  -0.25%  " },
    @{ Row = 5; D = "245.93"; E = "  -1.74%  " },
    @{ Row = 6; D = "0.661"; E = "  -1.30%  " },
    @{ Row = 7; D = "59.24"; E = "  -1.34%  " },
    @{ Row = 9; D = $null; E = "  -2.62%  " },
    @{ Row = 10; D = $null; E = "  -2.63%  " },
    @{ Row = 11; D = "0.110"; E = "  +2.15%  " },
    @{ Row = 12; D = "15.51"; E = $null },
    @{ Row = 13; D = "0.890"; E = "  +7.72%  " },
    @{ Row = 14; D = "2.358.40"; E = "  +0.07%  " },
    @{ Row = 15; D = $null; E = "  -0.76%  " },
    @{ Row = 16; D = "2.055.77"; E = "  -0.16%  " },
    @{ Row = 17; D = "18.20"; E = $null },
    @{ Row = 18; D = "37.004.57"; E = "  -0.36%  " },
    @{ Row = 19; D = "74.13"; E = "  -1.85%  " },
    @{ Row = 20; D = "0.0₃0892"; E = "  -2.13%  " },
    @{ Row = 21; D = $null; E = "  -0.14%  " },
    @{ Row = 22; D = "239.06"; E = "  +0.04%  " },
    @{ Row = 23; D = $null; E = "  +0.00%  " },
    @{ Row = 24; D = "2.46"; E = "  +1.72%  " },
    @{ Row = 25; D = "9.98"; E = "  +4.97%  " },
    @{ Row = 26; D = "169.28"; E = "  -0.04%  " },
    @{ Row = 27; D = $null; E = "  -4.34%  " },
    @{ Row = 28; D = "20.16"; E = "  -0.75%  " },
    @{ Row = 29; D = "5.44"; E = "  +11.80%  " },
    @{ Row = 30; D = $null; E = "  -1.45%  " },
    @{ Row = 31; D = "1.12"; E = "  -2.70%  " },
    @{ Row = 32; D = $null; E = "  +2.93%  " },
    @{ Row = 33; D = $null; E = "  -1.28%  " },
    @{ Row = 34; D = $null; E = "  +2.89%  " },
    @{ Row = 35; D = $null; E = "  +0.12%  " },
    @{ Row = 36; D = "1.84"; E = "  +5.93%  " },
    @{ Row = 37; D = "0.0839"; E = "  -6.16%  " },
    @{ Row = 38; D = $null; E = "  -0.91%  " },
    @{ Row = 39; D = $null; E = "  +1.69%  " },
    @{ Row = 40; D = "3.08"; E = "  -1.15%  " },
    @{ Row = 41; D = "0.0224"; E = "  -0.49%  " },
    @{ Row = 42; D = "1.17"; E = "  +1.58%  " },
    @{ Row = 43; D = "0.0964"; E = "  -10.95%  " },
    @{ Row = 44; D = "97.97"; E = "  -0.05%  " },
    @{ Row = 45; D = "17.08"; E = "  -4.34%  " },
    @{ Row = 46; D = "1.302.74"; E = "  +0.69%  " },
    @{ Row = 47; D = $null; E = "  -4.34%  " },
    @{ Row = 49; D = "6.82"; E = "  -1.03%  " },
    @{ Row = 50; D = "2.245.76"; E = "  +0.28%  " },
    @{ Row = 51; D = "44.51"; E = "  +1.93%  " }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($row, 4)
        # Force text storage so numeric-looking strings (e.g. "245.93",
        # "0.110") are preserved exactly instead of being normalised into
        # floating point numbers by Excel's auto-detection.
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($row, 5).Value = $u.E
    }
}
